$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Re-sort the data range (two rows were appended unsorted at the bottom
# - "MERSEN" and "SCHNEIDER ELECTRIC" - this sorts them alphabetically
# by column A, same as the rest of the list).
$sortRange = $ws.Range("A2:G51")
$keyRange = $ws.Range("A2:A51")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange) | Out-Null

$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = -4142   # xlNo
$ws.Sort.Orientation = 1  # xlSortRows (top to bottom)
$ws.Sort.Apply()

# Restore view state (scroll position / selection) similar to commit.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B35").Select() | Out-Null
